# Auto-generated edit script applying scheduled-runner value updates
# to the Halicarnassus_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 333.41666
$ws.Range("I80").Value = 266.4
$ws.Range("K80").Value = 799.1999999999999
$ws.Range("M80").Value = 198.8000000000001

# Row 83
$ws.Range("H83").Value = 333.41666
$ws.Range("I83").Value = 266.4
$ws.Range("K83").Value = 2397.6
$ws.Range("M83").Value = 2594.4

# Row 88
$ws.Range("H88").Value = 2791.8572
$ws.Range("I88").Value = 3151.6
$ws.Range("J88").Value = 1892.5
$ws.Range("K88").Value = 3151.6
$ws.Range("L88").Value = 1892.5
$ws.Range("M88").Value = -2745.6
$ws.Range("N88").Value = -2704.5

# Row 91
$ws.Range("H91").Value = 2791.8572
$ws.Range("I91").Value = 3151.6
$ws.Range("J91").Value = 1892.5
$ws.Range("K91").Value = 3151.6
$ws.Range("L91").Value = 1892.5
$ws.Range("M91").Value = -1747.6
$ws.Range("N91").Value = -4700.5

# Row 94
$ws.Range("H94").Value = 4712.154
$ws.Range("I94").Value = 3841.6365
$ws.Range("K94").Value = 3841.6365
$ws.Range("M94").Value = -3390.6365

# Row 98
$ws.Range("H98").Value = 530.6667
$ws.Range("I98").Value = 530.6667
$ws.Range("K98").Value = 530.6667
$ws.Range("M98").Value = 967.3333

# Row 107
$ws.Range("H107").Value = 260.55554
$ws.Range("I107").Value = 258.125
$ws.Range("J107").Value = 280
$ws.Range("K107").Value = 258.125
$ws.Range("L107").Value = 280
$ws.Range("M107").Value = 1661.875
$ws.Range("N107").Value = -4120

# Row 116
$ws.Range("H116").Value = 3995
$ws.Range("I116").Value = 3993.3333
$ws.Range("K116").Value = 3993.3333
$ws.Range("M116").Value = -551.3332999999998

# Row 122
$ws.Range("H122").Value = 530.6667
$ws.Range("I122").Value = 530.6667
$ws.Range("K122").Value = 1592.0001
$ws.Range("M122").Value = 857.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 1076
$ws.Range("I97").Value = 1021.1
$ws.Range("J97").Value = 1350.5
$ws.Range("K97").Value = 1021.1
$ws.Range("L97").Value = 1350.5
$ws.Range("M97").Value = -525.1
$ws.Range("N97").Value = -2342.5

# Row 110
$ws.Range("H110").Value = 142860820
$ws.Range("I110").Value = 200002850
$ws.Range("J110").Value = 5750
$ws.Range("K110").Value = 200002850
$ws.Range("L110").Value = 5750
$ws.Range("M110").Value = -200000805
$ws.Range("N110").Value = -9840

$ws = $wb.Worksheets.Item("BSM")
# Row 12
$ws.Range("H12").Value = 433.33334
$ws.Range("I12").Value = 433.33334
$ws.Range("K12").Value = 433.33334
$ws.Range("M12").Value = -265.33334

# Row 94
$ws.Range("H94").Value = 1000
$ws.Range("I94").Value = 1000
$ws.Range("K94").Value = 1000
$ws.Range("M94").Value = -549

# Row 99
$ws.Range("H99").Value = 71429690
$ws.Range("I99").Value = 83334450
$ws.Range("K99").Value = 83334450
$ws.Range("M99").Value = -83332952

# Row 105
$ws.Range("H105").Value = 7396676
$ws.Range("I105").Value = 13866643
$ws.Range("J105").Value = 2428.5715
$ws.Range("K105").Value = 13866643
$ws.Range("L105").Value = 2428.5715
$ws.Range("M105").Value = -13864896
$ws.Range("N105").Value = -5922.5715

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 500
$ws.Range("K6").Value = 500
$ws.Range("M6").Value = -387

# Row 7
$ws.Range("H7").Value = 62.68182
$ws.Range("I7").Value = 39.714287
$ws.Range("J7").Value = 102.875
$ws.Range("K7").Value = 39.714287
$ws.Range("L7").Value = 102.875
$ws.Range("M7").Value = 73.285713
$ws.Range("N7").Value = -328.875

# Row 22
$ws.Range("H22").Value = 1932.6666
$ws.Range("I22").Value = 399
$ws.Range("K22").Value = 399
$ws.Range("M22").Value = -49

# Row 55
$ws.Range("H55").Value = 19131.5
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 19131.5
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 19131.5
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -19761.5

# Row 58
$ws.Range("H58").Value = 1354.2941
$ws.Range("I58").Value = 732.9655
$ws.Range("K58").Value = 732.9655
$ws.Range("M58").Value = -529.9655

# Row 68
$ws.Range("H68").Value = 86235
$ws.Range("J68").Value = 86235
$ws.Range("L68").Value = 86235
$ws.Range("N68").Value = -87733

# Row 71
$ws.Range("H71").Value = 86235
$ws.Range("J71").Value = 86235
$ws.Range("L71").Value = 258705
$ws.Range("N71").Value = -266193

# Row 134
$ws.Range("H134").Value = 2958
$ws.Range("I134").Value = 2119.889
$ws.Range("K134").Value = 6359.667
$ws.Range("M134").Value = -3824.667

# Row 136
$ws.Range("H136").Value = 1354.2941
$ws.Range("I136").Value = 732.9655
$ws.Range("K136").Value = 2198.8965
$ws.Range("M136").Value = 351.1035000000002

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 235029.11
$ws.Range("I4").Value = 251907.75
$ws.Range("K4").Value = 755723.25
$ws.Range("M4").Value = -755611.25

# Row 11
$ws.Range("H11").Value = 8026.3335
$ws.Range("I11").Value = 9477.772000000001
$ws.Range("K11").Value = 28433.316
$ws.Range("M11").Value = -28293.316

# Row 68
$ws.Range("H68").Value = 667.6667
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

# Row 71
$ws.Range("H71").Value = 667.6667
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

# Row 107
$ws.Range("H107").Value = 630.0714
$ws.Range("I107").Value = 324.4
$ws.Range("K107").Value = 973.1999999999999
$ws.Range("M107").Value = 946.8000000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 49
$ws.Range("H49").Value = 26398
$ws.Range("J49").Value = 26398
$ws.Range("L49").Value = 26398
$ws.Range("N49").Value = -26766

# Row 113
$ws.Range("H113").Value = 8877.777
$ws.Range("J113").Value = 8877.777
$ws.Range("L113").Value = 8877.777
$ws.Range("N113").Value = -13217.777

# Row 126
$ws.Range("H126").Value = 3170.7
$ws.Range("I126").Value = 2944
$ws.Range("K126").Value = 8832
$ws.Range("M126").Value = -6362

# Row 132
$ws.Range("H132").Value = 43868.25
$ws.Range("I132").Value = 49729.715
$ws.Range("K132").Value = 149189.145
$ws.Range("M132").Value = -146659.145

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4929.7
$ws.Range("I40").Value = 4922
$ws.Range("K40").Value = 4922
$ws.Range("M40").Value = -4786

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 400
$ws.Range("J2").Value = 400
$ws.Range("L2").Value = 400
$ws.Range("N2").Value = -624

# Row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

# Row 54
$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -21040

# Row 107
$ws.Range("H107").Value = 83334090
$ws.Range("I107").Value = 83334090
$ws.Range("K107").Value = 250002270
$ws.Range("M107").Value = -250000350

# Row 126
$ws.Range("H126").Value = 3354.4707
$ws.Range("I126").Value = 1456.909
$ws.Range("K126").Value = 4370.727000000001
$ws.Range("M126").Value = -1900.727000000001
